$d = $word.ActiveDocument

# --- (Technical definition) -> (Technical description), add spacing after=0 ---
$d.Content.Find.Execute("(Technical definition)", $false, $false, $false, $false, $false, $true, 1, $false, "(Technical description)", 2) | Out-Null
$d.Paragraphs(8).Range.ParagraphFormat.SpaceAfter = 0

# --- Intro paragraph: slidecast topic change ---
$d.Content.Find.Execute(
    "This memo contains information regarding the requested slidecast on the “Lock-Out, Tag-Out” methodology.  The next sections elaborate on the target audience and the purpose of the slidecast and on the rhetorical moves and organizing pattern used.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "This memo contains information regarding the requested slidecast on Nvidia GeForce RTX 3070 graphics card.  The next sections elaborate on the target audience and the purpose of the slidecast and on the rhetorical moves and organizing pattern used.",
    2) | Out-Null

# --- Audience paragraph ---
$d.Content.Find.Execute(
    "The target audience will be workers and managers in industrial manufacturing where hazardous energy is involved and poses a potential health risk to personnel.  The audience will already have an understanding of where they might encounter hazardous energies, how to perform their own duties safely, and how to perform specific LOTO procedures in their own workspace.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "The target audience will be consumers looking to purchase the latest Nvidia graphics card for personal or professional usage.  The audience will already have an understanding of the other components of a personal computer (PC) that are required to make a graphics processing unit (GPU) functional.",
    2) | Out-Null

# --- Purpose paragraph ---
$d.Content.Find.Execute(
    "My purpose is to write a technical definition for general “Lock-Out, Tag-Out” (LOTO) procedures in an industrial workplace.  I will include explanations as to what is LOTO, why LOTO is necessary for a functioning and safe work environment, and how LOTO can be implemented without disrupting workflow.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "My purpose is to provide a technical description of the RTX 3070 graphics card developed by Nvidia Corporation.  I will include specific details as to what can be included in a RTX 3070 graphics card, what a RTX 3070 could be used to do, where a consumer might be able to purchase one, and current price ranges.",
    2) | Out-Null

# --- Rhetorical moves intro ---
$d.Content.Find.Execute(
    "In this slide cast, I used the following rhetorical moves:",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "In this slide cast, I will use the following rhetorical moves:",
    2) | Out-Null

# --- Rhetorical move bullets ---
$d.Content.Find.Execute(
    "Explain the history of x.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Describe in detail the specific features of x.",
    2) | Out-Null

$d.Content.Find.Execute(
    "Explain pros/cons or advantages/disadvantages of x.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Describe in detail the specific characteristics of x.",
    2) | Out-Null

$d.Content.Find.Execute(
    "Explain something that is necessary to use x.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Describe in detail specific behaviors of x.",
    2) | Out-Null

$d.Content.Find.Execute(
    "Use an example of x.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Describe in detail the specifics of what x is made of.",
    2) | Out-Null

$d.Content.Find.Execute(
    "Compare x to something similar.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Use a metaphor to explain what x is like.",
    2) | Out-Null

# --- Organizing pattern paragraph ---
$d.Content.Find.Execute(
    "In the slidecast, I introduce the concept of “Lock-Out, Tag-Out” before expanding into its reasons and purposes.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "In the slidecast, I introduce the RTX 3070 before getting into greater details.",
    2) | Out-Null
